# Update "想去人数" (interested-attendee counts) in column F across sheets,
# reflecting a refreshed scrape of the source data (commit: "Update gh-pages
# to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 0
$ws1.Range("F3").Value = 101
$ws1.Range("F4").Value = 64
$ws1.Range("F5").Value = 0
$ws1.Range("F8").Value = 74
$ws1.Range("F9").Value = 0
$ws1.Range("F11").Value = 22
$ws1.Range("F14").Value = 0
$ws1.Range("F16").Value = 417
$ws1.Range("F18").Value = 41
$ws1.Range("F19").Value = 0
$ws1.Range("F20").Value = 5280
$ws1.Range("F21").Value = 122
$ws1.Range("F22").Value = 178
$ws1.Range("F23").Value = 694
$ws1.Range("F24").Value = 220
$ws1.Range("F25").Value = 246

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 46

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7000
$ws4.Range("F3").Value = 101
$ws4.Range("F5").Value = 457
$ws4.Range("F6").Value = 157
$ws4.Range("F7").Value = 6911
$ws4.Range("F8").Value = 74
$ws4.Range("F9").Value = 0
$ws4.Range("F12").Value = 0
$ws4.Range("F14").Value = 150
$ws4.Range("F15").Value = 0
$ws4.Range("F16").Value = 0
$ws4.Range("F17").Value = 50
$ws4.Range("F19").Value = 0
$ws4.Range("F21").Value = 5280
$ws4.Range("F22").Value = 0
$ws4.Range("F23").Value = 122
$ws4.Range("F24").Value = 178
$ws4.Range("F25").Value = 0
$ws4.Range("F26").Value = 220
$ws4.Range("F27").Value = 246
